$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PoiFormulaHelperTest")

# Update formula in A6: (34*45)+(235/65) -> (34*45)+(800/40)
$ws.Range("A6").Formula = "=(34*45)+(800/40)"

# Move the active selection from A7 to A6
$ws.Activate()
$ws.Range("A6").Select()
